$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style of A5 (the last existing row label cell) onto A6 so the
# new "Ensemble" label matches the formatting of the other model names.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A6").Value = "Ensemble"
$ws.Range("B6").Value = 0.3
$ws.Range("C6").Value = -0.01
$ws.Range("D6").Value = 0.444
$ws.Range("E6").Value = 0.666
$ws.Range("F6").Value = 0.672
$ws.Range("G6").Value = 0.553
